$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.768.45"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "1.634.46"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.14"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("E6").Value = "  -0.55%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("E8").Value = "  +0.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0638"
$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.68"
$ws.Range("E10").Value = "  -3.55%  "

$ws.Range("E11").Value = "  +1.27%  "

$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "1.635.05"
$ws.Range("E13").Value = "  -0.70%  "

$ws.Range("D14").Value = "1.859.01"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.86"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").Value = "25.780.73"

$ws.Range("E19").Value = "  -0.25%  "

$ws.Range("E20").Value = "  +1.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.48"
$ws.Range("E21").Value = "  +0.80%  "

$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.28"
$ws.Range("E23").Value = "  +2.23%  "

$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("E25").Value = "  +3.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.77"

$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.89"
$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.60"
$ws.Range("E29").Value = "  +0.57%  "

$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.26"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.59"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.904"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").Value = "1.128.67"
$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("E39").Value = "  -1.64%  "

$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("E42").Value = "  +2.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.26"
$ws.Range("E43").Value = "  +1.02%  "

$ws.Range("E44").Value = "  +0.75%  "

$ws.Range("D45").Value = "1.768.75"
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("D46").Value = "0.0₆0109"
$ws.Range("E46").Value = "  -2.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.19"
$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("E48").Value = "  -2.37%  "

$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.59"
$ws.Range("E50").Value = "  -2.87%  "

$ws.Range("E51").Value = "  +2.26%  "
